$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    9997.246906068069,
    10173.45744822817,
    10339.93953956769,
    10431.95027711086,
    10592.37669578395,
    10746.00605383846,
    10271.87855223233,
    10369.10025149762,
    10532.14637421962,
    10407.58267755104,
    10632.81133896388,
    10651.1099562725,
    10726.85278796381,
    10828.75712841147,
    10882.98531550991,
    10890.62910734961,
    10942.57470522964,
    10991.56801295198,
    11023.61004017605,
    11067.69167583976,
    11130.51835252798,
    11188.64494197969,
    11244.32732485725,
    11326.26999923725,
    11317.27012166006,
    11352.86163011717,
    11390.24358065011,
    11423.33321721212,
    11446.53281407104,
    11471.75498974751,
    11504.39738237219,
    11532.65579059415,
    11553.73698554174,
    11546.035698489,
    11541.61392115616,
    11572.08635561047,
    11559.62338037264,
    11586.16365151377,
    11593.47519524045,
    11649.39243731134,
    11626.46082571328,
    11634.53273713433,
    11651.48835333896,
    11655.79836002462,
    11669.77365811561,
    11714.87602733246,
    11700.85087783734,
    11714.23676616257,
    11730.35249669618,
    11753.82501792888,
    11767.72335051302,
    11790.28777418434,
    11806.90118486327,
    11829.74660756819,
    11856.38309275827,
    11851.24584566233,
    11870.78751915871,
    11887.16992177463
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

Write-Output "Updated column C rows 2-59 with new transit depth values."
